$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.354337453842163
$ws.Range("B1").Value = 1.552693247795105
$ws.Range("C1").Value = 1.948180913925171
$ws.Range("D1").Value = 1.994009613990784
$ws.Range("E1").Value = 1.642878651618958
